# Build site at 2022-09-26 16:07:08 UTC
# This script reproduces the edit applied to LOM3016.xlsx:
#   - The three "Docentes responsaveis" rows (old rows 13-15, which only had
#     values in columns B/C) are removed.
#   - As a consequence, everything below shifts up by three rows.
#   - The "Objetivos:" row (row 10) content is replaced by the first
#     professor's name, and the professor names / evaluation text are
#     rearranged into the rows that now follow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the three rows that used to hold the "Docentes responsaveis" names
#    (old rows 13, 14, 15). Deleting the same row index three times removes
#    all three because each deletion shifts the following rows up.
$ws.Range("A13:C13").EntireRow.Delete()
$ws.Range("A13:C13").EntireRow.Delete()
$ws.Range("A13:C13").EntireRow.Delete()

# 2) Update the text content so the final layout matches the target state.
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("B13").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C13").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("B15").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C15").Value = "5983729 - Fernando Vernilli Junior"

$ws.Range("B18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("B19").Value = "Serão aplicadas duas provas escritas com notas P1 e P2."
$ws.Range("C19").Value = "Serão aplicadas duas provas escritas com notas P1 e P2."

$ws.Range("B20").Value = "A nota final NF será calculada pela fórmula: NF=(P1 + P2)/2."
$ws.Range("C20").Value = "A nota final NF será calculada pela fórmula: NF=(P1 + P2)/2."

$ws.Range("B21").Value = "Será aplicada uma prova escrita NR que comporá com a nota final NF a média final após recuperação MF=(NF+NF)/2."
$ws.Range("C21").Value = "Será aplicada uma prova escrita NR que comporá com a nota final NF a média final após recuperação MF=(NF+NF)/2."
